$d = $word.ActiveDocument

# Locate the three affected list-item paragraphs by their current text so the
# script is resilient to exact paragraph-index shifts. They sit consecutively:
#   "interaction"                   -> "Title screen / saving"
#   "Title screen / saving"         -> "decision making"
#   "Dialogue and decision making"  -> (paragraph removed entirely)

$paras = $d.Paragraphs
$count = $paras.Count

$idxInteraction = -1
$idxTitleScreen = -1
$idxDialogue = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $paras.Item($i).Range.Text.TrimEnd("`r", "`n", [char]7)
    if ($t -eq "interaction") {
        $idxInteraction = $i
    } elseif ($t -eq "Title screen / saving" -and $idxInteraction -ge 0 -and $idxTitleScreen -lt 0) {
        $idxTitleScreen = $i
    } elseif ($t -eq "Dialogue and decision making" -and $idxTitleScreen -ge 0 -and $idxDialogue -lt 0) {
        $idxDialogue = $i
    }
}

if ($idxInteraction -gt 0) {
    $paras.Item($idxInteraction).Range.Text = "Title screen / saving"
}

if ($idxTitleScreen -gt 0) {
    $paras.Item($idxTitleScreen).Range.Text = "decision making"
}

if ($idxDialogue -gt 0) {
    $paras.Item($idxDialogue).Range.Delete()
}
